$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix / verify idTrajet (column A) values ---
$ws.Range("A2").Value = 1001
$ws.Range("A3").Value = 1002
$ws.Range("A4").Value = 0

# --- New columns: Etat du trajet (L) / Etat du chauffeur (M) ---
$ws.Range("L1").Value = "Etat du trajet"
$ws.Range("M1").Value = "Etat du chauffeur"

$ws.Range("L2").Value = """Pas commencé"""
$ws.Range("M2").Value = 0

$ws.Range("L3").Value = """En cours"""
$ws.Range("M3").Value = 1

$ws.Range("L4").Value = """Finis"""
$ws.Range("M4").Value = 1

# --- Column widths for the newly inserted columns ---
# (COM ColumnWidth snaps to whole-pixel increments, so these are the closest
# achievable values to the authored 15.85546875 / 18.140625 widths.)
$ws.Columns.Item(12).ColumnWidth = 14.91667
$ws.Columns.Item(13).ColumnWidth = 17.25

# --- Update the selected cell shown when the workbook is reopened ---
$ws.Range("E23").Select()
